# "feedback form and other component updated"
# Replace the single student record (row 2) with a new student's data,
# clear out the old second record (row 3, except formatting on E3),
# repoint the E2 mailto hyperlink to the new student's e-mail address,
# and update the view (scroll position / selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Student")

# ---- Row 2: new student record --------------------------------------
$ws.Range("A2").Value2 = "16EGICS039"
$ws.Range("B2").Value2 = "Himanshu Panchal"
$ws.Range("C2").Value2 = "Mr. Rahul Panchal"
$ws.Range("D2").Value2 = "Mrs. Jaya Panchal"
$ws.Range("E2").Value2 = "panchalhimanshu@gmail.com"
$ws.Range("F2").Value2 = 7410258963
$ws.Range("G2").Value2 = 9874563210
$ws.Range("H2").Value2 = 8520369147
$ws.Range("I2").Value2 = 111111
$ws.Range("J2").Value2 = 222222
$ws.Range("K2").Value2 = "sdjfhbhjb"

# ---- Row 3: remove the old second record, keep E3 formatting --------
$ws.Range("A3:D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("F3:K3").ClearContents()

# ---- Hyperlinks: only E2 should remain, pointing at the new e-mail --
$ws.Hyperlinks.Delete() | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:panchalhimanshu@gmail.com") | Out-Null
$ws.Range("E2").Style = "Hyperlink"

# ---- View: scroll so column E is at the left, select L2 -------------
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$ws.Range("L2").Select() | Out-Null
